# edit communication uart and edit file exel frame
#
# The sheet holds a "UART frame" reference table (rows 1-7). This edit:
#   - tweaks a few labels in that table (END BYTE / 1 control / DLEN = 5 /
#     clears the old "0 check" note),
#   - moves the lone formatted spacer cell that used to live at M13 up to
#     M12 (row 13 gets reused for the new table below),
#   - adds a second, condensed copy of the frame table in rows 13-18
#     (dropping the DATA column, relabelled "DLEN = 3"),
#   - leaves the active selection on D4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Create the two new merged areas up front (done first so the style they
# incidentally touch doesn't get mixed into the real content styles below).
# ---------------------------------------------------------------------
$ws.Range("A13:B13").Merge()
$ws.Range("C18:E18").Merge()

# ---------------------------------------------------------------------
# Row 13 <- header row 1 (HEAD BYTE / DLEN / DirDATA / CRC8bit / END BYTE)
# ---------------------------------------------------------------------
$ws.Range("A1").Copy($ws.Range("A13"))
$ws.Range("B1").Copy($ws.Range("B13"))
$ws.Range("C1").Copy($ws.Range("C13"))
$ws.Range("D1").Copy($ws.Range("D13"))
$ws.Range("G1").Copy($ws.Range("E13"))

# H1's label changes to "END BYTE" - do that now so the copy into F13
# below (and the rest of the script) picks up the new text.
$ws.Range("H1").Value = "END BYTE"
$ws.Range("H1").Copy($ws.Range("F13"))

# ---------------------------------------------------------------------
# Row 14 <- data row 2 (0x55 / 1byte / 1byte / 1byte / 0xFF)
# ---------------------------------------------------------------------
$ws.Range("A2").Copy($ws.Range("A14"))
$ws.Range("B2").Copy($ws.Range("B14"))
$ws.Range("C2").Copy($ws.Range("C14"))
$ws.Range("D2").Copy($ws.Range("D14"))
$ws.Range("G2").Copy($ws.Range("E14"))
$ws.Range("H2").Copy($ws.Range("F14"))

# ---------------------------------------------------------------------
# Row 15 <- blank spacer (formatting only)
# ---------------------------------------------------------------------
$ws.Range("D5").Copy($ws.Range("D15"))
$ws.Range("H3").Copy($ws.Range("F15"))

# ---------------------------------------------------------------------
# Row 16 <- "0 check" label (copy D4 before it is cleared below)
# ---------------------------------------------------------------------
$ws.Range("D4").Copy($ws.Range("D16"))
$ws.Range("H4").Copy($ws.Range("F16"))

# ---------------------------------------------------------------------
# Row 17 <- blank spacer (formatting only)
# ---------------------------------------------------------------------
$ws.Range("D5").Copy($ws.Range("D17"))
$ws.Range("H5").Copy($ws.Range("F17"))

# ---------------------------------------------------------------------
# Row 18 <- "DLEN = 3" banner: same fill as the "DLEN" header cell (C1),
# centered, across the merged C18:E18 area.
# ---------------------------------------------------------------------
$ws.Range("C1").Copy($ws.Range("C18:E18"))
$ws.Range("C18:E18").HorizontalAlignment = -4108
$ws.Range("C18").Value = "DLEN = 3"

# ---------------------------------------------------------------------
# Remaining label edits on the first table.
# ---------------------------------------------------------------------
$ws.Range("D3").Value = "1 control"
$ws.Range("D4").ClearContents()
$ws.Range("C6").Value = "DLEN = 5"

# ---------------------------------------------------------------------
# Move the formatted spacer cell from M13 to M12, then drop the old one.
# ---------------------------------------------------------------------
$ws.Range("M12").Orientation = 0
$ws.Range("M13").Clear()

# ---------------------------------------------------------------------
# Selection ends on D4.
# ---------------------------------------------------------------------
$ws.Range("D4").Select()
